$wb = $excel.ActiveWorkbook

# --- Sheet: Details ---
$wsDetails = $wb.Worksheets.Item("Details")

# Row 3
$wsDetails.Range("W3").Value = "10.1.2.0 - "

# Row 5
$wsDetails.Range("U5").Value = "yes"
$wsDetails.Range("W5").Value = "10.1.2.0 - "

# Row 6
$wsDetails.Range("M6").Value = "6-8"
$wsDetails.Range("N6").Value = "30"
$wsDetails.Range("O6").Value = "yes"
$wsDetails.Range("P6").Value = "yes"
$wsDetails.Range("Q6").Value = "3"
$wsDetails.Range("R6").Value = "yes"
$wsDetails.Range("U6").Value = "yes"
$wsDetails.Range("W6").Value = "10.1.2.0 - "
$wsDetails.Range("X6").Value = "yes"
$wsDetails.Range("Y6").Value = "lead2=admin,lead3=user,Pradeep=user,karthik=admin,lead1=user,lead1=user,test1=admin"
$wsDetails.Range("Z6").Value = "2"

# New unique shared-string values, entered in the same order the
# original author typed them (this fixes their position in the
# shared strings table): a,b ; A ; 8-6 ; 2-4
$wsDetails.Range("AA6").Value = "a,b"
$wsDetails.Range("AA5").Value = "A"
$wsDetails.Range("L7").Value = "8-6"
$wsDetails.Range("L6").Value = "2-4"

$wsDetails.Activate()
$wsDetails.Range("G12").Select()

# --- Sheet: ResourceName ---
$wsResource = $wb.Worksheets.Item("ResourceName")
$wsResource.Activate()
$wsResource.Range("B12").Select()
